# Ind_finan_Saldos_Nuble.xlsx — add 27 monthly rows (2018-10 .. 2020-12) above the
# existing data block, shifting the current rows (2021-01 .. 2024-08) down by 27.
#
# The sheet "Cuadro" currently has data rows 4..47 (2021-01-01 .. 2024-08-01).
# After the edit it must have data rows 4..74 (2018-10-01 .. 2024-08-01), i.e. the
# 44 existing rows move to 31..74 and 27 brand-new rows are written into 4..30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 4
$lastDataRow  = 47
$shift = 27

# --- 1) Shift the existing data rows down by $shift, bottom-up so sources are
#        read before they get overwritten. Range copy/paste carries the cell
#        style along (s="5" for the date column, s="11" for the numeric ones),
#        which is exactly what the new destination rows need too. ---
for ($r = $lastDataRow; $r -ge $firstDataRow; $r--) {
    $dst = $r + $shift
    $srcRange = "A" + $r + ":F" + $r
    $dstRange = "A" + $dst + ":F" + $dst
    $ws.Range($srcRange).Copy()
    $ws.Range($dstRange).PasteSpecial(-4104)
}

# --- 2) New rows' numeric data (columns B..F), taken straight from the diff. ---
$newRows = @(
    @(2018,10,46160,2081704.4756499,125,18102591.432,50809.176811),
    @(2018,11,46428,2015776.2572154,128,16110170.117187,52044.987722),
    @(2018,12,46783,2038928.4409721,133,14497128.097744,57024.704748),
    @(2019,1,47132,2070439.013494,135,13797354.948148,55954.101131),
    @(2019,2,47358,2285232.7359474,135,14729263.148148,57493.174753),
    @(2019,3,47828,2273053.4139416,136,14950762.625,59737.281988),
    @(2019,4,48246,2201518.9556439,138,9650000.4927536,59519.057701),
    @(2019,5,48736,2164976.8780778,140,9745533.3714285,63849.324772),
    @(2019,6,49092,2169277.8624419,140,9185208.0357142,65237.894262),
    @(2019,7,49559,2118374.9948546,140,10494164.935714,61135.451145),
    @(2019,8,49988,2071974.9345643,139,10818102.366906,64123.758469),
    @(2019,9,50600,2116498.4579051,145,8838295.7241379,60909.456888),
    @(2019,10,50946,2116380.0662858,147,8558201.0680272,63968.128026999),
    @(2019,11,51143,2122997.4600238,149,9279552.7248322,64612.681257),
    @(2019,12,51314,2183966.6062673,150,10644746.6,69356.825935),
    @(2020,1,51587,2190126.9910636,152,11867948.815789,71705.741782),
    @(2020,2,51733,2366502.5614404,157,13351531.076433,73210.853954),
    @(2020,3,51869,2387998.6100561,163,11097064.490797,75735.997657),
    @(2020,4,51808,2496531.4760847,166,10738512.825301,86888.035258),
    @(2020,5,51823,2742874.7586013,172,16290071.377907,91999.653427),
    @(2020,6,51893,2971798.1986202,173,15296816.046242,90305.722351),
    @(2020,7,52066,3058891.0203203,177,14549090.101694,99231.89643),
    @(2020,8,52256,3214355.4516419,181,15332235,172930.8846),
    @(2020,9,52445,3521764.5205262,184,13722888.706521,174024.247263),
    @(2020,10,52773,3623745.5599075,190,11905771.668421,160882.197197),
    @(2020,11,53121,3567764.5514768,194,12667853.113402,146903.297939),
    @(2020,12,53481,3653013.1572521,200,12310696.645,213092.804615)
)

$row = $firstDataRow
foreach ($rec in $newRows) {
    $ws.Cells.Item($row, 2).Value = $rec[2]
    $ws.Cells.Item($row, 3).Value = $rec[3]
    $ws.Cells.Item($row, 4).Value = $rec[4]
    $ws.Cells.Item($row, 5).Value = $rec[5]
    $ws.Cells.Item($row, 6).Value = $rec[6]
    $row = $row + 1
}

# --- 3) Column A (Periodo) for every data row, old and new alike, is rewritten
#        as a real Excel date (first of each month) anchored at 2018-10-01 in
#        row 4. The loader cannot round-trip the workbook's original ISO-8601
#        "t=d" date cells, so every row's date is restored explicitly here. ---
$lastRowNow = $lastDataRow + $shift
$year = 2018
$month = 10
for ($r = $firstDataRow; $r -le $lastRowNow; $r++) {
    $d = Get-Date -Year $year -Month $month -Day 1 -Hour 0 -Minute 0 -Second 0
    $ws.Cells.Item($r, 1).Value = $d
    $month = $month + 1
    if ($month -gt 12) {
        $month = 1
        $year = $year + 1
    }
}
